$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text relabel ("Item" -> "item", "quantity" -> "quantity_measure") ---
$ws.Range("A1").Value = "item"
$ws.Range("B1").Value = "quantity_measure"

# --- C1 adopts the same bold/centered header style already used by A1 & B1 ---
$ws.Range("C1").Font.Italic = $false
$ws.Range("C1").HorizontalAlignment = -4108   # xlCenter

# --- D1 / E1 no longer carry the header fill/style - drop them entirely ---
$ws.Range("D1:E1").Clear()

# --- Add a thin spacer row (row 5) below the data, and extend the used range
#     down to it (touching a cell is what registers the new sheet extent) ---
$ws.Cells.Item(5, 1).Font.Bold = $false
$ws.Rows.Item(5).RowHeight = 4.2

# --- Update the active selection to B2:B4 ---
$ws.Range("B2:B4").Select()
